$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Convertidor")

# H4: was "3/8"" (text) -> now a plain number 8
$ws.Range("H4").Value = 8
# J4: 0.15 -> 0.605
$ws.Range("J4").Value = 0.60499999999999998

# J5: 0.6 -> 0.3
$ws.Range("J5").Value = 0.3

# H6: was "1/2"" (text) -> now "8mm" (text)
$ws.Range("H6").Value = "8mm"
# J6: 0.4 -> 0.6
$ws.Range("J6").Value = 0.6

# H14: was "3/8"" (text) -> now a plain number 8
$ws.Range("H14").Value = 8
# J14: 0.15 -> 0.605
$ws.Range("J14").Value = 0.60499999999999998
